# "add file log, agrreate table"
# The dim_products data-dictionary table (Table1, A10:D56) on sheet "dim_products"
# had its "VARCHAR" data-type entries changed to "TEXT".  After this edit the
# shared string "VARCHAR" is no longer referenced anywhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dim_products")

# Column C ("Data type") rows that previously read "VARCHAR" now read "TEXT".
$ranges = @("C11:C12", "C14:C22", "C24", "C26:C31", "C33:C53", "C55")
foreach ($r in $ranges) {
    $ws.Range($r).Value = "TEXT"
}

# Restore the active-cell position recorded in the sheet view.
$ws.Range("D16").Select()
